$d = $word.ActiveDocument

# --- Change 1a: insert new bullet paragraph "Instead of words, images emojis
#     can be sent as an easier way to communicate." after "Emojis generate
#     emotions in a digital conversation." (new 4th bullet under "advantages").
$pEmotions = $d.Paragraphs(20)
$pEmotions.Range.InsertParagraphAfter()
$pInstead = $d.Paragraphs(21)
$pInstead.Range.Text = "Instead of words, images emojis can be sent as an easier way to communicate."

# --- Change 1b: add a new leading run "Minimize input required from user."
#     inside the paragraph that currently only holds a bold <w:br/> run
#     (this paragraph moved down to index 22 after the insertion above).
$pBreak = $d.Paragraphs(22)
$insertStart = $pBreak.Range.Start
$minimizeText = "Minimize input required from user."
$pBreak.Range.InsertBefore($minimizeText)
$newRunRange = $d.Range($insertStart, $insertStart + $minimizeText.Length)
$newRunRange.Font.Bold = 0

# --- Change 2: insert new bullet paragraph "Emojis don’t expresses specific
#     words." after "Technical glitches can change the appearance of emojis..."
#     (new 5th bullet under "disadvantages"). One paragraph was inserted
#     above this point (the "Minimize..." edit only added a run, not a new
#     paragraph), so the target paragraph shifted from 27 to 28.
$pTechnical = $d.Paragraphs(28)
$pTechnical.Range.InsertParagraphAfter()
$pDontExpress = $d.Paragraphs(29)
$pDontExpress.Range.Text = "Emojis don’t expresses specific words."

# --- Change 3 & 4: move the _GoBack bookmark from the very last paragraph of
#     the document to the empty bullet paragraph that sits right after the
#     two "S" bullets (before "List and explain three (3) examples of how
#     emoji equity affects users."). Two more paragraphs were inserted above
#     this point, so it shifted from 36 to 38.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$pGoBackTarget = $d.Paragraphs(38)
$d.Bookmarks.Add("_GoBack", $pGoBackTarget.Range)
